# Apply "added result of billing report LPH_V2.1.10B" edit:
#  - Column H ("Env"/"LPH") becomes a "Remark" column with its per-row
#    values cleared out (the LPH billing-report pass no longer stamps an
#    environment tag in every row).
#  - The active-cell selection moves from F16 to G15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Env" -> "Remark"
$ws.Cells.Item(1, 8).Value = "Remark"

# Rows 2-12: clear the "LPH" marker in column H, but keep it looking like the
# plain (non-bold, non-filled) bordered cells used elsewhere in the table —
# copy the format from the neighbouring plain cell (F) and then blank it out.
for ($r = 2; $r -le 12; $r++) {
    $plainCell = $ws.Cells.Item($r, 6)
    $hCell = $ws.Cells.Item($r, 8)
    $plainCell.Copy($hCell)
    $hCell.ClearContents()
}

# Move the saved selection / active cell
[void]$ws.Range("G15").Select()
